$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.056.68'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '2.645.39'
$ws.Range('E3').Value = '  +4.04%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '518.70'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.65'
$ws.Range('E6').Value = '  +1.32%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.566'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').Value = '2.660.27'
$ws.Range('E9').Value = '  +4.36%  '
$ws.Range('E10').Value = '  +2.79%  '
$ws.Range('E11').Value = '  +2.96%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.337'
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('E13').Value = '  -1.60%  '
$ws.Range('D14').Value = '3.124.11'
$ws.Range('E14').Value = '  +4.57%  '
$ws.Range('D15').Value = '59.207.22'
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.99'
$ws.Range('E16').Value = '  +1.31%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000137'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('D18').Value = '2.657.48'
$ws.Range('E18').Value = '  +4.51%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '349.97'
$ws.Range('E19').Value = '  +2.80%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.52'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.36'
$ws.Range('E21').Value = '  +2.46%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.18'
$ws.Range('E22').Value = '  +3.70%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '62.18'
$ws.Range('E24').Value = '  +2.39%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.420'
$ws.Range('E25').Value = '  +2.06%  '
$ws.Range('D26').Value = '2.758.45'
$ws.Range('E26').Value = '  +3.97%  '
$ws.Range('E27').Value = '  +1.81%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('D29').Value = '0.0₃0807'
$ws.Range('E29').Value = '  +2.45%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.15'
$ws.Range('E30').Value = '  +2.63%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.27'
$ws.Range('E32').Value = '  +7.40%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.02'
$ws.Range('E33').Value = '  +2.61%  '
$ws.Range('E34').Value = '  +2.83%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '149.99'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.962'
$ws.Range('E36').Value = '  +5.35%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.02'
$ws.Range('E37').Value = '  +3.28%  '
$ws.Range('E38').Value = '  +2.29%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.75'
$ws.Range('E39').Value = '  +1.75%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.843'
$ws.Range('E40').Value = '  +2.44%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.71'
$ws.Range('E41').Value = '  +5.27%  '
$ws.Range('E42').Value = '  +1.74%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.996'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '276.68'
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.609'
$ws.Range('E45').Value = '  +1.36%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0984'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '19.65'
$ws.Range('E47').Value = '  +5.09%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0523'
$ws.Range('E48').Value = '  -1.75%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0230'
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '10.29'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.994.11'
$ws.Range('E51').Value = '  +4.41%  '
